$wb = $excel.ActiveWorkbook

# --- 1) "Horas insumidas": log 3 new hours for S-01020 / "Métrica AUX_TM" ---
$wsHoras = $wb.Worksheets.Item("Horas insumidas")
$wsHoras.Activate()

$wsHoras.Range("B99").Value = 40478
$wsHoras.Range("B99").NumberFormat = "d-mmm"
$wsHoras.Range("C99").Value = "Duilio"
$wsHoras.Range("D99").Value = "Métrica AUX_TM"
$wsHoras.Range("E99").Value = "S-01020"
$wsHoras.Range("F99").Value = 3

$wsHoras.Range("D99").Select()
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1

# --- 2) "Earned Value": S-01020 task is now complete (100%) ---
$wsEV = $wb.Worksheets.Item("Earned Value")
$wsEV.Activate()

$wsEV.Range("C19").Value = "Completada"
$wsEV.Range("D19").Value = 100

$wsEV.Range("D19").Select()
